{"js": "// IMC 2018 Shadow PC reviews\n// 1) Bold the \"Paper4: S6, NSDI'18:\" label and append the review text as a\n//    separate (non-bold) run in the same paragraph.\n// 2) Add a new \"Paper 10: ENVI, ANCS'18: \" heading (preceded by a blank\n//    spacer paragraph) right after the Split/Merge paragraph.\n\nconst body = context.document.body;\n\n// ---- Change 1: Paper4 label + review text -------------------------------\nconst paper4Results = body.search(\"Paper4: S6, NSDI\\u201918: \", { matchCase: true });\npaper4Results.load(\"items\");\nawait context.sync();\n\nif (paper4Results.items.length > 0) {\n  const paper4Range = paper4Results.items[0];\n\n  // Strip the trailing space from the label and make it bold.\n  const labelRange = paper4Range.insertText(\"Paper4: S6, NSDI\\u201918:\", \"Replace\");\n  labelRange.font.bold = true;\n  await context.sync();\n\n  // Append the review text (leading space restored) as a non-bold run.\n  const bodyText =\n    \" Extensoin to OpenNF, Split/Merge, E2 and StatelessNF. This one argues \" +\n    \"that the prior methods suffer from performance (StatelessNF), high \" +\n    \"downtimes during scaling (OpenNF and Split/Merge) and limited NF \" +\n    \"functionality support (E2). The paper proposes a distributed shared \" +\n    \"state approach where the state is distributed among NF instances and \" +\n    \"resides in a global name, encapsulated in objects. Any NF can access \" +\n    \"and modify the the state objects as the objects are in local address \" +\n    \"space. Not sure how different it is from StatelessNF. \";\n  const reviewRange = labelRange.insertText(bodyText, \"After\");\n  reviewRange.font.bold = false;\n  await context.sync();\n}\n\n// ---- Change 2: new \"Paper 10\" heading after the Split/Merge paragraph ---\nconst splitMergeResults = body.search(\"Split/Merge has four components\", {\n  matchCase: true\n});\nsplitMergeResults.load(\"items\");\nawait context.sync();\n\nif (splitMergeResults.items.length > 0) {\n  const splitMergePara = splitMergeResults.items[0].paragraphs.getFirst();\n  await context.sync();\n\n  // Blank spacer paragraph (matches the style already used in the doc).\n  const blankPara = splitMergePara.insertParagraph(\"\", \"After\");\n  await context.sync();\n\n  // New heading paragraph: bold title run + trailing plain-space run.\n  const paper10Para = blankPara.insertParagraph(\n    \"Paper 10: ENVI, ANCS\\u201918:\",\n    \"After\"\n  );\n  await context.sync();\n\n  const titleRange = paper10Para.getRange();\n  titleRange.font.bold = true;\n  await context.sync();\n\n  const trailingSpace = titleRange.insertText(\" \", \"After\");\n  trailingSpace.font.bold = false;\n  await context.sync();\n}\n", "ps1": "# IMC 2018 Shadow PC reviews\n# 1) Bold the \"Paper4: S6, NSDI'18:\" label and append the review text as a\n#    separate (non-bold) run in the same paragraph.\n# 2) Add a new \"Paper 10: ENVI, ANCS'18: \" heading (preceded by a blank\n#    spacer paragraph) right after the Split/Merge paragraph.\n\n$d = $word.ActiveDocument\n\n# ---- Change 1: Paper4 label + review text --------------------------------\n$find = $d.Content\n$found = $find.Find.Execute(\"Paper4: S6, NSDI\u201918: \")\nif ($found) {\n    # Strip the trailing space from the label and make it bold.\n    $find.Text = \"Paper4: S6, NSDI\u201918:\"\n    $find.Font.Bold = 1\n\n    $reviewText = \" Extensoin to OpenNF, Split/Merge, E2 and StatelessNF. This one argues that the prior methods suffer from performance (StatelessNF), high downtimes during scaling (OpenNF and Split/Merge) and limited NF functionality support (E2). The paper proposes a distributed shared state approach where the state is distributed among NF instances and resides in a global name, encapsulated in objects. Any NF can access and modify the the state objects as the objects are in local address space. Not sure how different it is from StatelessNF. \"\n\n    $find.Collapse(0)\n    $find.InsertAfter($reviewText)\n    $find.Font.Bold = 0\n}\n\n# ---- Change 2: new \"Paper 10\" heading after the Split/Merge paragraph ----\n$count = $d.Paragraphs.Count\n$targetIdx = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs($i).Range.Text\n    if ($t -like \"*Split/Merge has four components*\") {\n        $targetIdx = $i\n        break\n    }\n}\n\nif ($targetIdx -gt 0) {\n    # Blank spacer paragraph (matches the style already used in the doc).\n    $splitMergeRange = $d.Paragraphs($targetIdx).Range\n    $splitMergeRange.Collapse(0)\n    $splitMergeRange.InsertParagraphAfter()\n\n    # New heading paragraph, inserted right after the spacer.\n    $blankRange = $d.Paragraphs($targetIdx + 1).Range\n    $blankRange.Collapse(0)\n    $blankRange.InsertParagraphAfter()\n\n    # Title run: bold, no trailing space.\n    $newRange = $d.Paragraphs($targetIdx + 2).Range\n    $newRange.InsertBefore(\"Paper 10: ENVI, ANCS\u201918:\")\n\n    $titleRange = $d.Paragraphs($targetIdx + 2).Range\n    $titleRange.MoveEnd(1, -1) | Out-Null\n    $titleRange.Font.Bold = 1\n\n    # Trailing plain-space run, kept non-bold.\n    $titleRange.Collapse(0)\n    $titleRange.InsertAfter(\" \")\n    $spaceRange = $titleRange.Duplicate\n    $spaceRange.Font.Bold = 0\n}\n"}
